{"js": "// Update the date line and the 25 division problems in the practice table.\n// Target cells by (row, column) position so that duplicate problem text\n// (e.g. \"14\u00f75=\" appears twice) is handled unambiguously.\n\n// 1) Update the date/weekday line at the top of the document.\nconst dateResults = context.document.body.search(\"2025-11-22 Saturday\", { matchCase: true });\ndateResults.load(\"items\");\nawait context.sync();\nif (dateResults.items.length > 0) {\n  dateResults.items[0].insertText(\"2025-11-23 Sunday\", \"Replace\");\n}\n\n// 2) Update the division problems inside the table, addressed by cell\n//    position so duplicate text values are handled correctly.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Map of [rowIndex, columnIndex] -> new text, matching the diff exactly.\nconst cellUpdates = [\n  [0, 0, \"84\u00f73=\"],\n  [0, 1, \"20\u00f73=\"],\n  [0, 2, \"14\u00f74=\"],\n  [0, 3, \"46\u00f74=\"],\n  [0, 4, \"56\u00f76=\"],\n  [4, 0, \"32\u00f72=\"],\n  [4, 1, \"43\u00f79=\"],\n  [4, 2, \"36\u00f72=\"],\n  [4, 3, \"52\u00f76=\"],\n  [4, 4, \"96\u00f73=\"],\n  [8, 0, \"79\u00f75=\"],\n  [8, 1, \"90\u00f76=\"],\n  [8, 2, \"59\u00f73=\"],\n  [8, 3, \"45\u00f79=\"],\n  [8, 4, \"32\u00f72=\"],\n  [12, 0, \"19\u00f76=\"],\n  [12, 1, \"56\u00f75=\"],\n  [12, 2, \"99\u00f74=\"],\n  [12, 3, \"33\u00f79=\"],\n  [12, 4, \"40\u00f75=\"],\n  [16, 0, \"57\u00f79=\"],\n  [16, 1, \"34\u00f77=\"],\n  [16, 2, \"39\u00f73=\"],\n  [16, 3, \"44\u00f77=\"],\n  [16, 4, \"91\u00f77=\"],\n];\n\nfor (const [rowIndex, colIndex, newText] of cellUpdates) {\n  const cell = table.getCell(rowIndex, colIndex);\n  const range = cell.body.getRange();\n  range.insertText(newText, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Update the date line and the 25 division problems in the practice table.\n# Table cells are addressed by (row, column) position (1-based, Word COM\n# style) so duplicate problem text (e.g. \"14\u00f75=\" appears twice) is handled\n# unambiguously.\n\n$d = $word.ActiveDocument\n\n# 1) Update the date/weekday line at the top of the document.\n$find = $d.Content.Find\n$find.Text = \"2025-11-22 Saturday\"\n$find.Replacement.Text = \"2025-11-23 Sunday\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n# 2) Update the division problems inside the table, addressed by cell\n#    position so duplicate text values are handled correctly.\n$table = $d.Tables.Item(1)\n\n# Rows (1-based) that hold problem text; 3 blank rows follow each one.\n$dataRows = @(1, 5, 9, 13, 17)\n\n$newValues = @(\n    @(\"84\u00f73=\", \"20\u00f73=\", \"14\u00f74=\", \"46\u00f74=\", \"56\u00f76=\"),\n    @(\"32\u00f72=\", \"43\u00f79=\", \"36\u00f72=\", \"52\u00f76=\", \"96\u00f73=\"),\n    @(\"79\u00f75=\", \"90\u00f76=\", \"59\u00f73=\", \"45\u00f79=\", \"32\u00f72=\"),\n    @(\"19\u00f76=\", \"56\u00f75=\", \"99\u00f74=\", \"33\u00f79=\", \"40\u00f75=\"),\n    @(\"57\u00f79=\", \"34\u00f77=\", \"39\u00f73=\", \"44\u00f77=\", \"91\u00f77=\")\n)\n\nfor ($i = 0; $i -lt $dataRows.Length; $i++) {\n    $rowIndex = $dataRows[$i]\n    $rowValues = $newValues[$i]\n    for ($col = 1; $col -le 5; $col++) {\n        $cell = $table.Cell($rowIndex, $col)\n        $cell.Range.Text = $rowValues[$col - 1]\n    }\n}\n\nWrite-Output \"done\"\n"}
